# CM20_TestData_ManageBankStatements_21C.xlsx - "Add files via upload" edit
#
# The real-world diff is dominated by incidental re-save noise produced by a
# different Excel build (fileVersion/rupBuild, absPath, revisionPtr, window
# geometry, sub-pixel column width / row height drift on Input_Value and
# Output_Value, etc.) - none of that is reachable through the Excel object
# model, it simply falls out of which Excel binary happened to resave the
# file. The one deliberate, content-level change a user actually made is on
# the "Input_Value" sheet: the sample URL / username / password that used to
# sit in S2:U2 (with a live hyperlink on S2) were cleared out, and the
# hyperlink was removed - dropping 3 shared strings that are no longer
# referenced anywhere (https://edrx.fa.us2.oraclecloud.com,
# IBM_IMPLEMENTATION_USER, Oracle1234).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate() | Out-Null

# Select the range being cleared (mirrors the saved selection in the diff).
$ws.Range("S2:U2").Select() | Out-Null

# Remove the hyperlink that lived on S2 (Target: https://edrx.fa.us2.oraclecloud.com/).
$ws.Range("S2").Hyperlinks.Delete() | Out-Null

# Clear the URL / username / password sample values out of S2:U2.
$ws.Range("S2:U2").ClearContents() | Out-Null
